# Anonymise company name / client manager name in the report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Filter:" description cell (B2) — swap the real company /
#    vendor name and the real contractor/manager name for generic placeholders.
$ws.Range("B2").Value = "Vendor/Company Name = 'COMPANY NAME' AND First Name Is Not Blank AND Last Name Is Not Blank AND Contractor/Vendor Number Begins With 'CLIENT MANAGER' AND Event Date >= '20210524 4:19:34 AM' AND Event Date <= '20210525 12:19:34 PM'"

# 2) Walk every data row (the table body starts at row 8 and the used range
#    extends to row 144) and anonymise the "Vendor/Company Name" (col E) and
#    "Parade Text" (col G) values. Column H ("Contractor/Vendor Number")
#    keeps the real-looking "Stuart Brace" / "STUART BRACE" text untouched —
#    only the company name fields are anonymised.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 8; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "Company Name"
    $ws.Cells.Item($r, 7).Value = "Comp Name – C Mgr"
}
